$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 333, shifting existing rows 333:437 down to 334:438
$ws.Rows.Item(333).Insert()

# Populate the newly inserted row 333 with the new record
$ws.Range("A333").Value = 4
$ws.Range("B333").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C333").Value = "Los Lagos"
$ws.Range("D333").Value = 44985
$ws.Range("E333").Value = 10
$ws.Range("F333").Value = 100112045
$ws.Range("G333").Value = "Zapallo"
$ws.Range("H333").Value = "Paine"
$ws.Range("I333").Value = "1a (cosecha)"
$ws.Range("J333").Value = 1200
$ws.Range("K333").Value = 500
$ws.Range("L333").Value = 600
$ws.Range("M333").Value = 550
$ws.Range("N333").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O333").Value = "Región de O'Higgins"
$ws.Range("P333").Value = 550
$ws.Range("Q333").Value = 1
$ws.Range("R333").Value = "Hortaliza"
